{"js": "// Diary-style document: a new entry \"\u4eca\u5929\u5929\u6c14\u975e\u5e38\u597d\uff01\uff01\" is inserted as its\n// own paragraph right before the (former) last paragraph, and the former\n// last paragraph's text is changed to \"\u4eca\u5929\u4e0d\u77e5\u9053\u8981\u5e72\u4ec0\u4e48\" (keeping the\n// trailing _GoBack bookmark in place).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nif (count === 0) {\n  throw new Error(\"Document has no paragraphs.\");\n}\n\n// The paragraph immediately before the current last paragraph (i.e. the\n// \"\u6674\uff1a\" line). Inserting the new paragraph right after it (rather than\n// \"Before\" the last paragraph) makes the new paragraph pick up that\n// paragraph's east-Asian run-property hint, matching the source formatting.\nconst anchorIndex = count >= 2 ? count - 2 : count - 1;\nconst anchorParagraph = paragraphs.items[anchorIndex];\n\nanchorParagraph.insertParagraph(\"\u4eca\u5929\u5929\u6c14\u975e\u5e38\u597d\uff01\uff01\", \"After\");\nawait context.sync();\n\n// Re-fetch paragraphs since indices shifted after the insertion, then\n// update the (new) last paragraph's text in place so its bookmarks and\n// other paragraph-level properties are preserved.\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = refreshedParagraphs.items[refreshedParagraphs.items.length - 1];\nlastParagraph.insertText(\"\u4eca\u5929\u4e0d\u77e5\u9053\u8981\u5e72\u4ec0\u4e48\", \"Replace\");\nawait context.sync();\n", "ps1": "# Diary-style document: a new entry \"\u4eca\u5929\u5929\u6c14\u975e\u5e38\u597d\uff01\uff01\" is inserted as its\n# own paragraph right before the (former) last paragraph, and the former\n# last paragraph's text is changed to \"\u4eca\u5929\u4e0d\u77e5\u9053\u8981\u5e72\u4ec0\u4e48\" (keeping the\n# trailing _GoBack bookmark in place).\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nif ($count -lt 1) {\n    throw \"Document has no paragraphs.\"\n}\n\n# The paragraph immediately before the current last paragraph (i.e. the\n# \"\u6674\uff1a\" line). Inserting the new paragraph right after it (rather than\n# before the last paragraph) makes the new paragraph pick up that\n# paragraph's east-Asian run-property hint, matching the source formatting.\n$anchorIndex = $count - 1\nif ($anchorIndex -lt 1) { $anchorIndex = $count }\n$anchorParagraph = $d.Paragraphs.Item($anchorIndex)\n\n$anchorRange = $anchorParagraph.Range\n$anchorRange.Collapse(0)  # wdCollapseEnd\n$anchorRange.InsertParagraphAfter()\n\n# The freshly inserted (still empty) paragraph now sits right after the\n# anchor paragraph; give it the new diary-entry text.\n$newParagraph = $d.Paragraphs.Item($anchorIndex + 1)\n$newParagraph.Range.Text = \"\u4eca\u5929\u5929\u6c14\u975e\u5e38\u597d\uff01\uff01\"\n\n# Update the (new) last paragraph's text in place, excluding its trailing\n# paragraph mark, so the bookmark that follows it is preserved.\n$newCount = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($newCount)\n$lastRange = $lastParagraph.Range\n$lastRange.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1\n$lastRange.Text = \"\u4eca\u5929\u4e0d\u77e5\u9053\u8981\u5e72\u4ec0\u4e48\"\n"}
